# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 16 (pushing the existing
# rows 16-36 down to 17-37), then populate the new row with the latest
# week's data for Melón - Tuna - Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16..36 down to 17..37 by inserting a new row at 16.
$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 44557
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 100112027
$ws.Range("G16").Value = "Melón"
$ws.Range("H16").Value = "Tuna"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 9500
$ws.Range("N16").Value = "$/caja 18 unidades"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 528
$ws.Range("Q16").Value = 18
$ws.Range("R16").Value = "Hortaliza"
